# ============================================================
# Edit: add "2022-Q3" quarter data
#  1) Insert a new row on the "总计" summary sheet for 2022-Q3
#     (pushing the existing quarterly rows down by one)
#  2) Insert a brand-new "2022-Q3" worksheet (positioned right
#     after "总计") with the per-fund holding detail for the
#     quarter
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------
# 1. "总计" (summary) sheet: insert new row 2 for 2022-Q3
# ------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# carry the existing formatting down onto the freshly inserted row
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 34
$summary.Cells.Item(2, 4).Value = 9.35

# re-number the "序号" (index) column for the rows that shifted down
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ------------------------------------------------------------
# 2. Brand-new "2022-Q3" worksheet, placed right after "总计"
# ------------------------------------------------------------
$template = $wb.Worksheets.Item(2)   # "2022-Q2" -- same layout, used only as a format donor

$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# header row (B1:H1) + index column (A2:A35) reuse the bold/bordered
# style already used by every other quarterly sheet
$template.Cells.Item(1, 1).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$template.Cells.Item(2, 1).Copy()
$q3.Range("A2:A35").PasteSpecial(-4122)

$q3.Range("B1:H1").Value = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")

$q3Data = @(
        @('001216','易方达新收益灵活配置混合 - A','39.46','77.61','7.50','2.9595',2),
        @('009812','易方达悦兴一年持有期混合A','61.64','20.74','2.01','1.2390',4),
        @('001217','易方达新收益灵活配置混合 - C','14.06','77.61','7.50','1.0545',2),
        @('160527','博时研究优选3年封闭运作灵活配置混合A','16.93','97.11','5.71','0.9667',7),
        @('001603','易方达安盈回报混合','25.59','32.33','3.41','0.8726',3),
        @('010536','泰康优势企业混合A','12.44','94.53','5.07','0.6307',10),
        @('257010','国联安小盘精选混合','8.50','74.70','5.78','0.4913',3),
        @('009813','易方达悦兴一年持有期混合C','17.96','20.74','2.01','0.3610',4),
        @('016950','鹏华睿投灵活配置混合C','4.12','83.97','3.82','0.1574',1),
        @('010537','泰康优势企业混合C','2.37','94.53','5.07','0.1202',10),
        @('006013','易方达鑫转招利混合A','5.40','22.51','1.69','0.0913',3),
        @('002222','嘉实新趋势灵活配置混合','4.64','26.64','1.01','0.0469',10),
        @('001688','嘉实新起点灵活配置混合A','5.33','24.70','0.87','0.0464',8),
        @('001755','嘉实新思路灵活配置混合','4.85','25.81','0.88','0.0427',9),
        @('006138','国联安价值优选股票','0.57','94.64','5.99','0.0341',5),
        @('160528','博时研究优选3年封闭运作灵活配置混合C','0.59','97.11','5.71','0.0337',7),
        @('000573','天弘通利混合','1.01','79.25','3.13','0.0316',7),
        @('006014','易方达鑫转招利混合C','1.65','22.51','1.69','0.0279',3),
        @('002367','国联安安稳灵活配置混合','0.57','47.79','3.65','0.0208',2),
        @('014627','财通多策略福瑞混合（LOF）C','1.14','61.59','1.81','0.0206',2),
        @('013774','易方达趋势优选混合A','0.65','84.55','2.90','0.0188',10),
        @('501028','财通多策略福瑞混合（LOF）A','0.88','61.59','1.81','0.0159',2),
        @('562900','易方达中证现代农业主题ETF','0.64','97.55','2.35','0.0150',8),
        @('009750','汇安价值蓝筹混合A','0.30','94.35','4.60','0.0138',8),
        @('009751','汇安价值蓝筹混合C','0.19','94.35','4.60','0.0087',8),
        @('010154','中加中证500指数增强C','0.51','94.15','1.43','0.0073',9),
        @('003242','创金合信量化发现灵活配置混合C','0.40','92.08','1.65','0.0066',2),
        @('010153','中加中证500指数增强A','0.44','94.15','1.43','0.0063',9),
        @('003241','创金合信量化发现灵活配置混合A','0.32','92.08','1.65','0.0053',2),
        @('013775','易方达趋势优选混合C','0.11','84.55','2.90','0.0032',10),
        @('008113','中泰中证500指数增强C','0.44','91.57','0.70','0.0031',5),
        @('008112','中泰中证500指数增强A','0.16','91.57','0.70','0.0011',5),
        @('002178','嘉实新起点灵活配置混合C','0.01','24.70','0.87','0.0001',8),
        @('005434','鹏华睿投灵活配置混合A','0.00','83.97','3.82',0,1)
    )

# numeric "序号" (row index) column, 0-based
for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $q3.Cells.Item($i + 2, 1).Value = $i
}

# fund-code / name / size / position / ratio / value columns are stored
# as literal text in the source data (keeps the trailing zeroes, e.g.
# "7.50" must not collapse to 7.5) -- temporarily force text format,
# write the values, then drop back to the default "Normal" style so no
# extra number-format survives on the cell
$textRange = $q3.Range("B2:G35")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
}

$textRange.Style = "Normal"

# the smallest holding rounds to 0.0000 and is stored as a real number
# (0) rather than the text "0.0000" -- matches the convention used by
# every other quarterly sheet in this workbook
$q3.Cells.Item(35, 7).Value = 0

$q3.Range("A1").Select()
